$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Entities")
$ws2 = $wb.Worksheets.Item("Enum")

# --- Sheet "Entities" ---

# Quest reward now goes toward story progression: rewardKey01 (column R)
# gets -1 for both quest rows.
$ws1.Range("R2").Value = -1
$ws1.Range("R3").Value = -1

# Update the explanatory notes in column Y.
# Y2 (row2 note) - now documents the new reward-key convention.
$ws1.Range("Y2").Value = "rewardKey < 0 add towards story progression. Empty rewardKey gives nothing."

# Y1 (header note) - merge DialogueResponse + taskRequiredAmount guidance.
$ws1.Range("Y1").Value = "taskActionType of DialogueResponse, you need to fill in dialogue id into taskObjectiveKey. Fill in the response answer in taskRequiredAmount."

# Y3 was an empty styled placeholder cell - clear it (value + formatting) entirely.
$ws1.Range("Y3").Clear()

# Row 4 only held the now-removed "rewardKey that is empty..." note - delete it.
$ws1.Rows.Item(4).Delete()

# Resize (bestFit-style) columns so their widths reflect the new content:
# J grew a lot (longer merged note), R/S/T/U/V/W widened slightly, Y widened a lot.
$ws1.Range("J1").ColumnWidth = 36.0
$ws1.Range("R1").ColumnWidth = 10.833333333333334
$ws1.Range("S1").ColumnWidth = 14.5
$ws1.Range("T1").ColumnWidth = 10.833333333333334
$ws1.Range("U1").ColumnWidth = 14.5
$ws1.Range("V1").ColumnWidth = 10.833333333333334
$ws1.Range("W1").ColumnWidth = 14.5
$ws1.Range("Y1").ColumnWidth = 117.16666666666667

# --- Sheet "Enum" ---
# Update its selection first so re-selecting on "Entities" afterwards leaves
# "Entities" as the active/tabSelected sheet (matches the authored file).
$ws2.Range("B6").Select()

# --- Back to sheet "Entities" ---
# Selection tidy-up to match the authored workbook state (also drops the
# stale topLeftCell="K1" scroll position from the sheet view).
$ws1.Range("I13").Select()
